$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 83 (shifts rows 83:91 down to 84:92),
# Excel extends the existing merged-cell pattern of the row above automatically.
$ws.Rows.Item(83).Insert()

# Fill in the new row 83 for the new item "رول اون ريكسونا25"
$ws.Cells.Item(83, 1).Value = 77
$ws.Range("C83").Value = "رول اون ريكسونا25"
$ws.Range("H83").Value = "6:0"
$ws.Range("N83").Value = "35.00"
$ws.Range("P83").Value = "35.0000"
$ws.Range("Q83").Value = "1:0"

# Renumber the "م" (index) column for the rows that were pushed down (old 77..83 -> now 84..90)
for ($i = 84; $i -le 90; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 6
}

# Update the grand total (now row 91) to include the new item's price
$ws.Range("P91").Value = 5080.8549999999996

# Update the generated timestamp in the footer row (now row 92)
$ws.Range("A92").Value = "Monday, 15 September, 2025 6:58 PM"
